$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A (everything shifts right by one)
$ws.Columns("A:A").Insert()

# Copy the header style (from what is now B1, the old A1 header) onto the
# new column's header/data cells so the new column matches the existing
# look (bold header row, etc.)
$ws.Range("B1").Copy()
$ws.Range("A1:A16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header for the new column
$ws.Range("A1").Value = "Fig Index"

# Fig Index values, one per data row, pairing each row with its distortion
# metric / figure.
$ws.Range("A2").Value = "10.1038:s43018-023-00513-2_zero_fig1"
$ws.Range("A3").Value = "10.1038:s43018-023-00513-2_zero_fig2"
$ws.Range("A4").Value = "10.1038:s43018-023-00524-z_zero_fig1"
$ws.Range("A5").Value = "10.1038:s43018-023-00563-6_zero_fig1"
$ws.Range("A6").Value = "10.1038:s43018-023-00573-4_zero_fig1"
$ws.Range("A7").Value = "10.1038:s43018-023-00594-z_zero_fig1"
$ws.Range("A8").Value = "10.1038:s43018-023-00594-z_zero_fig2"
$ws.Range("A9").Value = "10.1038:s43018-023-00610-2_zero_fig1"
$ws.Range("A10").Value = "10.1038:s43018-023-00610-2_zero_fig2"
$ws.Range("A11").Value = "10.1038:s43018-023-00610-2_zero_fig3"
$ws.Range("A12").Value = "10.1038:s43018-023-00635-7_zero_fig1"
$ws.Range("A13").Value = "10.1038:s43018-023-00635-7_zero_fig2"
$ws.Range("A14").Value = "10.1038:s43018-023-00635-7_zero_fig3"
$ws.Range("A15").Value = "10.1038:s43018-023-00635-7_zero_fig4"
$ws.Range("A16").Value = "10.1038:s43018-023-00635-7_zero_fig5"

# Move the active selection like in the saved file
$ws.Range("C5").Select()
